$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 275 ("T03"/days_fu) so the new
# "Thromboprophylaxis" derived-variable entry slots in alphabetically,
# shifting every row from 275..294 down to 276..295.
$ws.Rows.Item(275).Insert() | Out-Null

# Populate the newly inserted row with the new derived variable.
$ws.Range("A275").Value = "Rx23"
$ws.Range("B275").Value = "thromboprophy"
$ws.Range("C275").Value = "Treatments"
$ws.Range("D275").Value = "Thromboprophylaxis"
# Column E (Values) is intentionally left blank for this row.

# Grow Table1 so its range / autofilter cover the new last row.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E295")) | Out-Null

# Reflect the view state captured after the edit (scroll position, zoom,
# active selection) as closely as the object model allows.
$ws.Range("B264").Select() | Out-Null
$excel.ActiveWindow.Zoom = 200
